# Correcting Relevance Markers Appenzeller-Herzog (2019) - van Dis (2020)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value = 0.525

$ws.Range("H3").Value = 0.1734317343173432
$ws.Range("I3").Value = 0.2798701298701299

$ws.Range("K3").Value = 85.15000000000001

$ws.Range("Q3").Value = 6
$ws.Range("R3").Value = 11
$ws.Range("S3").Value = 44
$ws.Range("T3").Value = 120

$ws.Range("V3").Value = 225
$ws.Range("W3").Value = 220
$ws.Range("X3").Value = 187
$ws.Range("Y3").Value = 111

$ws.Range("AF3").Value = 0.9740259999999999
$ws.Range("AG3").Value = 0.952381
$ws.Range("AH3").Value = 0.809524
$ws.Range("AI3").Value = 0.480519
